$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68
$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 'Conhecimentos Específicos'
$ws.Range("C68").Value = 'Logística'
$ws.Range("D68").Value = 'Conceito de Logística (<i>Council of Logistics Management</i>):'
$ws.Range("E68").Value = 'é o processo de controle, planejamento e implementação do fluxo'
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0

# Row 69
$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 'Conhecimentos Específicos'
$ws.Range("C69").Value = 'Logística'
$ws.Range("D69").Value = 'objetivo da logística, Ballou (2005):'
$ws.Range("E69").Value = 'dispor a mercadoria ou o serviço certo, no lugar certo, no tempo certo e nas condições desejadas'
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0

# Row 70
$ws.Range("A70").Value = 69
$ws.Range("B70").Value = 'Conhecimentos Específicos'
$ws.Range("C70").Value = 'Logística'
$ws.Range("D70").Value = 'tripé logístico'
$ws.Range("E70").Value = @'
<ul>
	<li>Transporte;</li>
	<li>Distribuição; e</li>
	<li>Armazenagem.</li>
</ul>
'@
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0

# Row 71
$ws.Range("A71").Value = 70
$ws.Range("B71").Value = 'Conhecimentos Específicos'
$ws.Range("C71").Value = 'Logística'
$ws.Range("D71").Value = 'divisão da logística segundo algumas literaturas'
$ws.Range("E71").Value = @'
<ul>
	<li>Atividades primárias: <ul> <li>Transportes,</li> <li>Processamento de Pedidos e</li> <li>Manutenção de estoques;</li> </ul> <i>Minemônico (TPM)</i></li>
	<li>Atividades de apoio: <ul> <li>Armazenagem,</li> <li>Manuseio de Materiais,</li> <li>Embalagem de Proteção,</li> <li>Obtenção, Programação de Produtos e</li> <li>Manutenção de Informação.</li> </ul></li>
</ul>
'@
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0

# Row 72
$ws.Range("A72").Value = 71
$ws.Range("B72").Value = 'Conhecimentos Específicos'
$ws.Range("C72").Value = 'Logística'
$ws.Range("D72").Value = 'Processo de fluxo de materias:'
$ws.Range("E72").Value = @'
<ol>
	<li>Entradas - Fornecedores</li>
	<li>Estoque/Arm azenamento</li>
	<li>Processo Produtivo</li>
	<li>Produtos acabados (depósito)</li>
	<li>Saídas - Clientes</li>
</ol>
'@
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0

# Row 73
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = 'Conhecimentos Específicos'
$ws.Range("C73").Value = 'Logística'
$ws.Range("D73").Value = 'Classe de materiais ao longo do processo produtivo'
$ws.Range("E73").Value = @'
<ul>
	<li>Matéria-prima</li>
	<li>Materiais em processamento</li>
	<li>Materiais semiacabados</li>
	<li>Materiais acabados</li>
	<li>Produtos acabados</li>
</ul>
'@
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0

# Row 74
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = 'Conhecimentos Específicos'
$ws.Range("C74").Value = 'Logística'
$ws.Range("D74").Value = @'
<b>Transporte:</b>
<i>Características</i>
'@
$ws.Range("E74").Value = @'
Conceito: <ul> <li>parte do processo logístico responsável por levar os produtos ao consumidor final e/ou entre fornecedor e produtor</li> </ul>
corresponde, em média, a 2/3 de todo o custo logístico
fatores determinantes em relação ao transporte é a escolha do modal: <ul> <li>Custo do transporte;</li> <li>Velocidade com que o produto é transportado;</li> <li>Tipo de manuseio do produto durante o transporte; e</li> <li>Quantidade de viagens</li> </ul>
'@
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0

# Row 75
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = 'Conhecimentos Específicos'
$ws.Range("C75").Value = 'Logística'
$ws.Range("D75").Value = 'Modais de Transporte:'
$ws.Range("E75").Value = @'
<ul> <li><b>Aeroviário </b><ul> <li>Vantages: <ul> <li>longas distâncias, independente da geografia</li> <li>mais rápido dentre os modais</li> <li>Menor custo com embalagens</li> </ul></li> <li>desvantagens: <ul> <li>volume pequeno</li> <li>custo mais elevado</li> <li>geralmente precisa de outro modal para concluir o transporte</li> </ul></li> </ul></li> <li><b>Aquaviário </b><ul> <li><b>Marítimo:</b> mares e oceanos</li> <li><b>Fluvial:</b> rios</li> <li><b>Lacustre:</b> lagos e lagoas.</li> <li>vantagens: <ul> <li>Maior capacidade de carga entre os modais;</li> <li>grandes distâncias de forma autônoma</li> <li>Baixo custo unitário de carregamento</li> </ul></li> <li>desvantagens: <ul> <li>mais lento entre os modais;</li> <li>Maior suscetibilidade as mudanças da natureza;</li> <li>Necessidade de terminais especializados</li> <li>Desembaraço burocrático</li> <li>Alto custo quanto ao seguro</li> </ul></li> </ul></li> <li><b>Ferroviário </b><ul> <li>vantagens: <ul> <li>Baixo custo</li> <li>Menor risco de acidentes</li> <li>grande capacidade</li> </ul></li> <li>desvantagens: <ul> <li>geralmente precisa de outro modal para concluir o transporte</li> <li>Baixo investimento governamental</li> <li>Rotas fixas e inflexíveis</li> </ul></li> </ul></li> <li><b>Rodoviário</b> <ul> <li>vantagens: <ul> <li>acessibilidade</li> <li>Rapidez para contratação</li> <li>Rotas flexíveis</li> <li>Menor burocracia entre os modais</li> <li>Menor custo de estrutura e alto investimento gorvenamental</li> </ul></li> <li>desvantagens: <ul> <li>Gastos com pedágio</li> <li>aumento de combustíveis tendem a aumentar o valor do frete</li> <li>capacidade de carga é bem menor</li> <li>baixa autonomia de deslocamento</li> <li>Maior chance de extravio</li> </ul></li> </ul></li> <li><b>Dutoviário </b><ul> <li>vantagens: <ul> <li>transporta grande volume de carga de forma constante</li> <li>grande confiabilidade no processo</li> <li>Baixos custos operacionais;</li> </ul></li> <li>Desvantagens: <ul> <li>Custo inicial para implantação altíssimo</li> <li>Burocracia ambiental</li> <li>Reduzida flexibilidade de trajeto.</li> </ul></li> </ul></li> </ul>

'@
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 1

# Row 76
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = '75'
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = 'Conhecimentos Específicos'
$ws.Range("C76").Value = 'Logística'
$ws.Range("D76").Value = 'Classificação dos Modais:'
$ws.Range("E76").Value = '<ul> <li><b>Velocidade de Transporte:</b> <ol> <li>Aeroviário</li> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aquaviário</li> <li>Dutoviário</li> </ol></li> <li><b>Disponibilidade:</b> <ol> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aeroviário</li> <li>Aquaviário</li> <li>Dutoviário</li> </ol></li> <li><b>Confiabilidade: </b><ol> <li>Dutoviário</li> <li>Rodoviário</li> <li>Ferroviário</li> <li>Aquaviário</li> <li>Aeroviário</li> </ol></li> <li><b>Capacidade de Carga:</b> <ol> <li>Aquaviário</li> <li>Ferroviário</li> <li>Rodoviário</li> <li>Aeroviário</li> <li>Dutoviário</li> </ol></li> <li><b>Frequência:</b> <ol> <li>Dutoviário</li> <li>Rodoviário</li> <li>Aeroviário</li> <li>Ferroviário</li> <li>Aquaviário</li> </ol></li> </ul>'
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 2

# Row 77
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = 'Conhecimentos Específicos'
$ws.Range("C77").Value = 'Logística'
$ws.Range("D77").Value = @'
<b>Distribuição</b>
<i>Características</i>
'@
$ws.Range("E77").Value = @'
<ul>
	<li><b>Caonceito:</b> <ul> <li>conjunto de ações voltadas à gestão de materiais, iniciando com a saída do produto do processo produtivo e terminando com a entrega no ponto final de consumo</li> </ul></li>
	<li><b>fatores mais importantes ligados à distribuição</b> <ul> <li>Conferência de cargas;</li> <li>Gestão do frete;</li> <li>Gestão do transporte;</li> <li>Análise e desempenho de indicadores;</li> <li>Gestão de Rotas ou Roteirização.</li> </ul></li>
</ul>

'@
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0

# Row 78
$ws.Range("A78").Value = 77
$ws.Range("B78").Value = 'Conhecimentos Específicos'
$ws.Range("C78").Value = 'Logística'
$ws.Range("D78").Value = @'
<b>Armazenamento</b>
<i>Características</i>
'@
$ws.Range("E78").Value = @'
<ul>
	<li>atividades que compreende a armazenagem: <ul> <li>receber</li> <li>carregar</li> <li>descarregar</li> <li>conservar</li> </ul></li>
	<li>quatro pontos principais para que uma empresa decida destinar uma parte de sua área útil à armazenagem, Ballou (1993): <ul> <li>reduzir custos de transporte e produção</li> <li>coordenação de suprimento e demanda</li> <li>auxílio ao processo de produção</li> <li><u>auxílio ao processo de marketing.</u></li> </ul></li>
	<li>funções da armazenagem, Ballou (1993): <ul> <li>Abrigo de produtos</li> <li>Consolidação</li> <li>Transferência e Transbordo</li> <li>Agrupamento</li> </ul></li>
	<li><b>codificação:</b> <ul> <li>catalogar, simplificar, especificar, normatizar e padronizar todo o estoque</li> <li>11 dígitos: <ol> <li>XX - Grupo</li> <li>XX - Classe</li> <li>XXXXXX - Código de identificação</li> <li>X - Dígito de Controle</li> </ol></li> </ul></li>
	<li><b>embalagens:</b> <ul> <li>vantagens: <ul> <li>proteção ao produto <ul> <li>manuseio</li> <li>transporte</li> <li>armazenagem</li> </ul></li> </ul></li> </ul></li>
	<li>ações pelos quais passam os materiais armazenados: <ul> <li>Especificação</li> <li>Simplificação</li> <li>Codificação</li> <li>Padronização</li> <li>Catalogação</li> <li>Normalização</li> </ul></li>
	<li><b>sistemas de armazenamento:</b> <ul> <li>Sistema WMS</li> <li>Racks</li> <li>Mezanino</li> <li>Sistema de carrossel</li> <li>Porta-paletes</li> <li>Flow Rack</li> </ul></li>
</ul>
'@
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0

# Row 79
$ws.Range("A79").Value = 78
$ws.Range("B79").Value = 'Conhecimentos Específicos'
$ws.Range("C79").Value = 'Logística'
$ws.Range("D79").Value = 'Sistemas Logísticos'
$ws.Range("E79").Value = @'
<ul>
	<li>TMS</li>
	<li>WMS</li>
	<li>Sistema de monitoramento de cargas</li>
	<li>Sistemas de roteirização</li>
	<li>Sistemas de gestão de frotas.</li>
</ul>
'@
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0

# Row 80
$ws.Range("A80").Value = 79
$ws.Range("B80").Value = 'Conhecimentos Específicos'
$ws.Range("C80").Value = 'Logística'
$ws.Range("D80").Value = @'
<b>Logística Reversa</b>
<i>Características</i>
'@
$ws.Range("E80").Value = @'
<ul>
	<li>responsabilidade sobre os resíduos produzidos em decorrência do consumo de bens</li>
	<li>devolução, reciclagem e adequada destinação de produtos pós-venda e pós consumo.</li>
	<li>etapas: <ol> <li>Devolução da embalagem ou resíduo para o comerciante;</li> <li>O comerciante devolve para ao fabricante; e</li> <li>O fabricante destina para reuso, reciclagem ou descarte adequado.</li> </ol></li>
	<li>leis que devem ser cumpridas</li>
	<li>preocupação com a lucratividade e sustentabilidade desse processo</li>
	<li>reversa: <ul> <li>transporte dos produtos nas mãos dos clientes de volta para a empresa</li> </ul></li>
	<li>atividades: <ul> <li>aterro sanitário</li> <li>doação</li> <li>processamento das devoluções</li> <li>reciclagem</li> <li>reembalagem</li> <li>remanufatura</li> <li>revenda</li> <li>revitalização</li> <li>recuperação de cargas roubadas ou perdidas</li> </ul></li>
</ul>
'@
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0

# Row 81
$ws.Range("A81").Value = 80
$ws.Range("B81").Value = 'Conhecimentos Específicos'
$ws.Range("C81").Value = 'Logística'
$ws.Range("D81").Value = @'
<b>Logística Verde</b>
<i>Conceito</i>
'@
$ws.Range("E81").Value = 'procedimentos de logística que objetivam a preservação do meio ambiente, que incluem desde a embalagem até o modal de transporte utilizado'
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0

# Row 82
$ws.Range("A82").Value = 81
$ws.Range("B82").Value = 'Conhecimentos Específicos'
$ws.Range("C82").Value = 'Logística'
$ws.Range("D82").Value = @'
<b>Logística de Pós-consumo</b>
<i>Características</i>
'@
$ws.Range("E82").Value = @'
<ul>
	<li>favorece o retorno dos produtos após serem utilizados pelos clientes, visando: <ul> <li>reciclagem</li> <li>reutilização ou</li> <li>descarte apropriado</li> </ul></li>
	<li><b>motivação:</b> <ul> <li>quantidade de materiais descartados pela sociedade desde o século XX até os dias de hoje <ul> <li>Diminuilção do ciclo de vida dos produtos</li> </ul> </li> </ul></li>
	<li>destinos para um produto após descarte: <ul> <li>local seguro (aterro sanitário)</li> <li>local não seguro</li> <li>Reciclagem</li> </ul></li>
</ul>
'@
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0

# Row 83
$ws.Range("A83").Value = 82
$ws.Range("B83").Value = 'Conhecimentos Específicos'
$ws.Range("C83").Value = 'Logística'
$ws.Range("D83").Value = @'
<b>Logística de Pós-venda</b>
<i>Características</i>
'@
$ws.Range("E83").Value = '<b>uma das suas preocupações-chave: </b><ul> <li>criar um canal acessível para clientes retornarem <b>produtos</b>. <ul> <li>defeitos de fabricação ou</li> <li>erros no pedido</li> </ul></li> </ul>'
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 0
